# Natmi following Dr Hou advice
# Recompute the Col8a1-Itga1 ligand/receptor table for the 3 clusters
# (ECs, FAPs, sCs) as sender AND target, producing all 9 sender/target
# combinations (rows 2-10) in place of the previous 6-row partial table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Col8a1", "Itga1", "ECs", 2, 0.6666666666666666, 4.573149333333333, 13.719448, 0.06199325535045235, 0.06199325535045234, 3, 1, 57.77686466666668, 173.330594, 0.5795494107546101, 0.57954941075461, 264.2222301324569, 2378.000071192112, 0.03592815460911474, 0.03592815460911473),
    @("ECs", "Col8a1", "Itga1", "FAPs", 2, 0.6666666666666666, 4.573149333333333, 13.719448, 0.06199325535045235, 0.06199325535045234, 3, 1, 9.278736333333333, 27.836209, 0.09307334701450438, 0.09307334701450438, 42.43304687695911, 381.897421892632, 0.005769919767791432, 0.005769919767791432),
    @("ECs", "Col8a1", "Itga1", "sCs", 2, 0.6666666666666666, 4.573149333333333, 13.719448, 0.06199325535045235, 0.06199325535045234, 3, 1, 32.637132, 97.91139600000001, 0.3273772422308855, 0.3273772422308855, 149.254478447712, 1343.290306029408, 0.02029518097354618, 0.02029518097354618),
    @("FAPs", "Col8a1", "Itga1", "ECs", 3, 1, 57.35861733333333, 172.075852, 0.7775489387534139, 0.7775489387534138, 3, 1, 57.77686466666668, 173.330594, 0.5795494107546101, 0.57954941075461, 3314.001071135121, 29826.00964021609, 0.4506280292874134, 0.4506280292874133),
    @("FAPs", "Col8a1", "Itga1", "FAPs", 3, 1, 57.35861733333333, 172.075852, 0.7775489387534139, 0.7775489387534138, 3, 1, 9.278736333333333, 27.836209, 0.09307334701450438, 0.09307334701450438, 532.2154866805631, 4789.939380125068, 0.0723690821973561, 0.0723690821973561),
    @("FAPs", "Col8a1", "Itga1", "sCs", 3, 1, 57.35861733333333, 172.075852, 0.7775489387534139, 0.7775489387534138, 3, 1, 32.637132, 97.91139600000001, 0.3273772422308855, 0.3273772422308855, 1872.020765245488, 16848.18688720939, 0.2545518272686444, 0.2545518272686443),
    @("sCs", "Col8a1", "Itga1", "ECs", 3, 1, 11.83673133333333, 35.510194, 0.1604578058961338, 0.1604578058961338, 3, 1, 57.77686466666668, 173.330594, 0.5795494107546101, 0.57954941075461, 683.889224341693, 6155.003019075237, 0.09299322685808194, 0.0929932268580819),
    @("sCs", "Col8a1", "Itga1", "FAPs", 3, 1, 11.83673133333333, 35.510194, 0.1604578058961338, 0.1604578058961338, 3, 1, 9.278736333333333, 27.836209, 0.09307334701450438, 0.09307334701450438, 109.8299090905051, 988.469181814546, 0.01493434504935685, 0.01493434504935684),
    @("sCs", "Col8a1", "Itga1", "sCs", 3, 1, 11.83673133333333, 35.510194, 0.1604578058961338, 0.1604578058961338, 3, 1, 32.637132, 97.91139600000001, 0.3273772422308855, 0.3273772422308855, 386.316962974536, 3476.852666770824, 0.052530233988695, 0.052530233988695)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowData = $rows[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
